$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "PoFDCtAE": fuel-balancing-priority updates. A batch of formulas that
# pulled percentages from 'Data from BFPIaE' are cleared back to 0 (hardcoded
# literal), and the dependent "1-x" helper cells are refreshed to match.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PoFDCtAE")
$ws.Activate()

# Formulas replaced outright with a literal 0
$ws.Range("C3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("S19").Value = 0
$ws.Range("T20").Value = 0

# These two keep their "1 - x" formula (their precedent now evaluates to 0,
# so they recompute to 1)
$ws.Range("R10").Formula = "=1-J10"
$ws.Range("R11").Formula = "=1-K11"

# These lose their formula entirely and become a hardcoded 0
$ws.Range("R14").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("R20").Value = 0

# Plain literal cell that was simply retyped from 1 to 0
$ws.Range("R18").Value = 0

# Selection left on A15 within the frozen bottom-right pane
$null = $ws.Range("A15").Select()

# ---------------------------------------------------------------------------
# Sheet "About": window was scrolled down (to row 49) before saving.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 1

# Restore PoFDCtAE as the active sheet/tab
$ws.Activate()
